$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.388.67"
$ws.Range("E2").Value = "  -0.86%  "

# Row 3
$ws.Range("D3").Value = "3.832.73"
$ws.Range("E3").Value = "  -2.48%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.24"
$ws.Range("E5").Value = "  +5.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.24"
$ws.Range("E6").Value = "  -6.45%  "

# Row 7
$ws.Range("E7").Value = "  -3.34%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("E9").Value = "  -5.12%  "

# Row 10
$ws.Range("E10").Value = "  -6.56%  "

# Row 11
$ws.Range("E11").Value = "  -8.48%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.07"
$ws.Range("E12").Value = "  -4.97%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.19"
$ws.Range("E13").Value = "  -2.93%  "

# Row 14
$ws.Range("D14").Value = "4.443.02"
$ws.Range("E14").Value = "  -2.47%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.50"
$ws.Range("E15").Value = "  +7.60%  "

# Row 16
$ws.Range("D16").Value = "3.834.60"
$ws.Range("E16").Value = "  -2.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.08"
$ws.Range("E17").Value = "  -1.67%  "

# Row 18
$ws.Range("E18").Value = "  -2.18%  "

# Row 19
$ws.Range("E19").Value = "  +3.24%  "

# Row 20
$ws.Range("D20").Value = "68.383.22"
$ws.Range("E20").Value = "  -0.89%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "414.03"
$ws.Range("E21").Value = "  -5.15%  "

# Row 22
$ws.Range("E22").Value = "  -2.77%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.84"
$ws.Range("E23").Value = "  -5.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.97"
$ws.Range("E24").Value = "  -3.93%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.89"
$ws.Range("E25").Value = "  +4.38%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.25"
$ws.Range("E26").Value = "  -8.48%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.40"
$ws.Range("E27").Value = "  -6.37%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.07"
$ws.Range("E28").Value = "  -5.67%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.01"
$ws.Range("E29").Value = "  -2.92%  "

# Row 30
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "669.15"
$ws.Range("E30").Value = "  -5.86%  "

# Row 31
$ws.Range("E31").Value = "  -6.37%  "

# Row 32
$ws.Range("E32").Value = "  -3.28%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.93"
$ws.Range("E33").Value = "  +6.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.21"
$ws.Range("E34").Value = "  +1.12%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.436"
$ws.Range("E35").Value = "  -9.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.30"
$ws.Range("E36").Value = "  -4.20%  "

# Row 37
$ws.Range("E37").Value = "  -8.92%  "

# Row 38
$ws.Range("B38").Value = "ThetaToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  +8.78%  "

# Row 39
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.15%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.146"
$ws.Range("E40").Value = "  -2.13%  "

# Row 41
$ws.Range("E41").Value = "  -0.07%  "

# Row 42
$ws.Range("E42").Value = "  -4.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.09"
$ws.Range("E43").Value = "  +2.60%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.74"
$ws.Range("E44").Value = "  -7.81%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.37"
$ws.Range("E45").Value = "  -1.44%  "

# Row 46
$ws.Range("E46").Value = "  -4.39%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.90"
$ws.Range("E47").Value = "  -2.97%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.22"
$ws.Range("E48").Value = "  +0.31%  "

# Row 49
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.23"
$ws.Range("E49").Value = "  -4.06%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.669.05"
$ws.Range("E50").Value = "  +9.79%  "

# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000260"
$ws.Range("E51").Value = "  +6.17%  "
